# The workbook's "Table1" (on sheet "Todos los Puntos de Medicion") is
# extended from 3 columns (Código, Orden, Ruta) to 5 columns:
#   Código, Orden, Definición de recorrido, Equipo, Descripción
# i.e. the "Ruta" column is renamed to "Definición de recorrido" and two
# brand-new columns "Equipo" and "Descripción" are appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$lo = $ws.ListObjects.Item(1)

# Rename the existing third column header ("Ruta" -> "Definición de recorrido").
$ws.Cells.Item(1, 3).Value = "Definición de recorrido"

# Append two new columns to the table; this grows the table/autofilter range
# from A1:C12 to A1:E12 automatically.
$colEquipo = $lo.ListColumns.Add()
$colDescripcion = $lo.ListColumns.Add()

# Name the newly added headers.
$ws.Cells.Item(1, 4).Value = "Equipo"
$ws.Cells.Item(1, 5).Value = "Descripción"

# Match the selection saved with the file (header row fully selected).
$ws.Activate()
$ws.Range("A1:E1").Select()
